# cs-en-us-062pct.xlsx weekly update
# - bump report volume/week-number and date range in the header
# - refresh the "Crime Complaints" table (rows 15-30 / Rape..Hate Crimes) with
#   the newly collected weekly figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Set a cell to a literal text value (e.g. "0" or "***.*") while keeping /
# adopting the "text" look (General format) taken from a sibling cell that
# already holds text, so the stored value is a genuine string and not a
# number that happens to render as "0".
function Set-TextCell($ws, $addr, $text, $fmtSrcAddr) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $src = $ws.Range($fmtSrcAddr)
    $src.Copy()
    $c.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Set a cell to a numeric value, adopting the number format from a sibling
# cell (used when a cell previously held placeholder text and now holds a
# real number).
function Set-NumCellWithFormat($ws, $addr, $value, $fmtSrcAddr) {
    $c = $ws.Range($addr)
    $c.Value = $value
    $src = $ws.Range($fmtSrcAddr)
    $src.Copy()
    $c.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Plain numeric update - cell already has the right format/type.
function Set-Num($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------------
# Header: volume/number and reporting week
# ---------------------------------------------------------------------------

$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextCell $ws "D15" "0" "C15"
Set-TextCell $ws "E15" "***.*" "C15"
Set-TextCell $ws "F15" "0" "C15"
Set-Num $ws "H15" -100
Set-Num $ws "M15" 45.454545454545

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-Num $ws "C16" 2
Set-Num $ws "D16" 3
Set-Num $ws "E16" -33.333333333333
Set-Num $ws "G16" 13
Set-Num $ws "H16" -23.076923076923
Set-Num $ws "I16" 106
Set-Num $ws "J16" 97
Set-Num $ws "K16" 9.278350515463
Set-Num $ws "L16" 12.765957446808
Set-Num $ws "M16" -31.612903225806
Set-Num $ws "N16" -84.342688330871

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-TextCell $ws "C17" "0" "C15"
Set-Num $ws "D17" 7
Set-Num $ws "E17" -100
Set-Num $ws "F17" 19
Set-Num $ws "G17" 18
Set-Num $ws "H17" 5.555555555555
Set-Num $ws "J17" 161
Set-Num $ws "K17" 2.484472049689
Set-Num $ws "L17" 32
Set-Num $ws "M17" 43.478260869565
Set-Num $ws "N17" -40.860215053763

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-Num $ws "C18" 5
Set-Num $ws "D18" 5
Set-Num $ws "F18" 17
Set-Num $ws "H18" 21.428571428571
Set-Num $ws "I18" 172
Set-Num $ws "J18" 129
Set-Num $ws "K18" 33.333333333333
Set-Num $ws "L18" 31.297709923664
Set-Num $ws "M18" -33.590733590733
Set-Num $ws "N18" -88.137931034482

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
Set-Num $ws "C19" 13
Set-Num $ws "D19" 17
Set-Num $ws "E19" -23.529411764705
Set-Num $ws "F19" 39
Set-Num $ws "G19" 47
Set-Num $ws "H19" -17.021276595744
Set-Num $ws "I19" 656
Set-Num $ws "J19" 506
Set-Num $ws "K19" 29.644268774703
Set-Num $ws "L19" 56.190476190476
Set-Num $ws "M19" 69.072164948453
Set-Num $ws "N19" -6.017191977077

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-Num $ws "C20" 1
Set-Num $ws "D20" 3
Set-Num $ws "E20" -66.666666666666
Set-Num $ws "F20" 10
Set-Num $ws "G20" 12
Set-Num $ws "H20" -16.666666666666
Set-Num $ws "I20" 135
Set-Num $ws "J20" 100
Set-Num $ws "K20" 35
Set-Num $ws "L20" 50
Set-Num $ws "M20" -5.594405594405
Set-Num $ws "N20" -91.646039603960

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
Set-Num $ws "C21" 21
Set-Num $ws "D21" 35
Set-Num $ws "E21" -40
Set-Num $ws "F21" 95
Set-Num $ws "G21" 107
Set-Num $ws "H21" -11.214953271028
Set-Num $ws "I21" 1255
Set-Num $ws "J21" 1015
Set-Num $ws "K21" 23.645320197044
Set-Num $ws "L21" 43.264840182648
Set-Num $ws "M21" 16.52739090065
Set-Num $ws "N21" -73.595623816536

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-NumCellWithFormat $ws "D22" 1 "G22"
Set-NumCellWithFormat $ws "E22" -100 "H22"
Set-Num $ws "G22" 4
Set-Num $ws "J22" 14
Set-Num $ws "K22" -21.428571428571
Set-Num $ws "M22" -57.692307692307

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
Set-Num $ws "C24" 45
Set-Num $ws "D24" 33
Set-Num $ws "E24" 36.363636363636
Set-Num $ws "F24" 173
Set-Num $ws "G24" 105
Set-Num $ws "H24" 64.761904761904
Set-Num $ws "I24" 1690
Set-Num $ws "J24" 1083
Set-Num $ws "K24" 56.048014773776
Set-Num $ws "L24" 57.063197026022
Set-Num $ws "M24" 89.461883408071

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
Set-Num $ws "C25" 10
Set-Num $ws "D25" 9
Set-Num $ws "E25" 11.111111111111
Set-Num $ws "F25" 42
Set-Num $ws "G25" 52
Set-Num $ws "H25" -19.230769230769
Set-Num $ws "I25" 414
Set-Num $ws "J25" 401
Set-Num $ws "K25" 3.241895261845
Set-Num $ws "L25" 38.461538461538
Set-Num $ws "M25" 4.810126582278

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextCell $ws "D26" "0" "C26"
Set-TextCell $ws "E26" "***.*" "C26"
Set-TextCell $ws "F26" "0" "C26"
Set-Num $ws "G26" 7
Set-Num $ws "H26" -100
Set-Num $ws "L26" -16

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-Num $ws "C27" 1
Set-Num $ws "D27" 1
Set-Num $ws "E27" 0
Set-Num $ws "F27" 6
Set-Num $ws "G27" 5
Set-Num $ws "H27" 20
Set-Num $ws "J27" 55
Set-Num $ws "K27" -12.727272727272

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-TextCell $ws "D30" "0" "C30"
Set-TextCell $ws "E30" "***.*" "C30"
